$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Extend the array formula in G9 so it spills into column H as well
# (G9:G12 -> G9:H12), which Excel back-fills with H9:H12 values.
$ws.Range("G9:H12").FormulaArray = "=G2:G5+E2:E5"

# New helper total in F11 that sums the G column's tail.
$ws.Range("F11").Formula = "=SUM(G10:G11,G12)"

# The original custom row heights (15.95) are no longer needed;
# auto-fitting drops the explicit height so the rows fall back to the
# sheet's default height.
$ws.Rows("2:5").AutoFit() | Out-Null
$ws.Rows("7:7").AutoFit() | Out-Null
$ws.Rows("9:9").AutoFit() | Out-Null

# Update the view: scroll/select so the active cell is G9 and the
# selection covers the newly extended array range.
$ws.Range("G9:H12").Select() | Out-Null
